$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 ("check that"), shifting existing rows down.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row with the "win conditions" task and its duration.
$ws.Range("B5").Value = "win conditions"
$ws.Range("C5").Value = 40

# Update the selected cell to match the target workbook.
$ws.Range("D8").Select()
